$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Rows.Item(9).Insert(-4121, 1)
$ws.Range("A10:F10").Copy()
$ws.Range("A9:F9").PasteSpecial(-4122)
Write-Host "done"
